$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries: "Butan" moves up (before "Timor Oriental"),
#     "Granada" shifts down into what was Butan's row.
#     Row 191 = Timor Oriental -> becomes Butan
#     Row 192 = Granada        -> becomes Timor Oriental
#     Row 193 = Butan          -> becomes Granada
$ws.Range("A191").Value = "Butan"
$ws.Range("B191").Value = 24
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 6
$ws.Range("E191").Value = 18

$ws.Range("A192").Value = "Timor Oriental"
$ws.Range("B192").Value = 24
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 24
$ws.Range("E192").Value = 0

$ws.Range("A193").Value = "Granada"
$ws.Range("B193").Value = 22
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 17
$ws.Range("E193").Value = 5

# --- Reorder countries: "Belice" moves up (before "Nueva Caledonia").
#     Row 199 = Nueva Caledonia -> becomes Belice
#     Row 200 = Belice          -> becomes Nueva Caledonia
$ws.Range("A199").Value = "Belice"
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2

$ws.Range("A200").Value = "Nueva Caledonia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# --- Update timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 06:35"

# --- Update country case numbers (India) ---
$ws.Range("B14").Value = 125149
$ws.Range("C14").Value = 355
$ws.Range("E14").Value = 69597
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 3728

# --- Update country case numbers (Pakistan) ---
$ws.Range("B22").Value = 52437
$ws.Range("C22").Value = 1743
$ws.Range("D22").Value = 16653
$ws.Range("E22").Value = 34683
$ws.Range("G22").Value = 34
$ws.Range("H22").Value = 1101

# --- Update country case numbers (Uzbekistan) ---
$ws.Range("B77").Value = 3036
$ws.Range("C77").Value = 8
$ws.Range("E77").Value = 531

# --- Update country case numbers (Camboya) ---
$ws.Range("B167").Value = 124
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 2

